# Update "想去人数" (number of interested attendees) figures for two
# events that appear on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8312   # 合肥·第十五届次元之门动漫游戏博览会: 8313 -> 8312
$wsExhibit.Range("F11").Value = 999   # 合肥·第九届环形宇宙动漫游戏嘉年华: 996 -> 999

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8312    # 合肥·第十五届次元之门动漫游戏博览会: 8313 -> 8312
$wsAll.Range("F15").Value = 999    # 合肥·第九届环形宇宙动漫游戏嘉年华: 996 -> 999
